$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B52: correct the imprecise floating point value to the rounded 1.52 ---
$ws.Range("B52").Value = 1.52

# --- Row 59 (year 2028): fill in the remaining monthly values C:M ---
$ws.Range("C59").Value = 2.18
$ws.Range("D59").Value = 2.61
$ws.Range("E59").Value = 5.65
$ws.Range("F59").Value = 10.47
$ws.Range("G59").Value = 13.22
$ws.Range("H59").Value = 14.03
$ws.Range("I59").Value = 13.52
$ws.Range("J59").Value = 10.26
$ws.Range("K59").Value = 7.26
$ws.Range("L59").Value = 10.11
$ws.Range("M59").Value = 1.12

# --- Row 60 (year 2029): fill in the monthly values B:M ---
$ws.Range("B60").Value = 1.87
$ws.Range("C60").Value = 2.53
$ws.Range("D60").Value = 0.97
$ws.Range("E60").Value = 5.29
$ws.Range("F60").Value = 10
$ws.Range("G60").Value = 10.99
$ws.Range("H60").Value = 10.93
$ws.Range("I60").Value = 13.26
$ws.Range("J60").Value = 12.15
$ws.Range("K60").Value = 11.15
$ws.Range("L60").Value = 6.85
$ws.Range("M60").Value = 1.3

# --- Row 61 (year 2030): fill in the monthly values B:M ---
$ws.Range("B61").Value = 2.77
$ws.Range("C61").Value = 3.36
$ws.Range("D61").Value = 2.41
$ws.Range("E61").Value = 2.95
$ws.Range("F61").Value = 8.800000000000001
$ws.Range("G61").Value = 12.33
$ws.Range("H61").Value = 12.35
$ws.Range("I61").Value = 14.54
$ws.Range("J61").Value = 10.39
$ws.Range("K61").Value = 8.42
$ws.Range("L61").Value = 5.71
$ws.Range("M61").Value = 3.87

# --- Row 62 (year 2031): fill in the monthly values B:E only ---
$ws.Range("B62").Value = 3.38
$ws.Range("C62").Value = 2.79
$ws.Range("D62").Value = 3.04
$ws.Range("E62").Value = 4.08
